$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44398
$ws.Range("J2").Value = 170
$ws.Range("K2").Value = 21000
$ws.Range("L2").Value = 22000
$ws.Range("M2").Value = 21500
$ws.Range("P2").Value = 538
$ws.Range("D3").Value = 44384
$ws.Range("H3").Value = "Madrigal"
$ws.Range("J3").Value = 80
$ws.Range("D4").Value = 44384
$ws.Range("H4").Value = "Madrigal"
$ws.Range("I4").Value = "Segunda"
$ws.Range("J4").Value = 30
$ws.Range("M4").Value = 19333
$ws.Range("P4").Value = 387
$ws.Range("D5").Value = 44384
$ws.Range("H5").Value = "Symphony"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 20000
$ws.Range("L5").Value = 21000
$ws.Range("M5").Value = 20400
$ws.Range("P5").Value = 510
$ws.Range("D6").Value = 44412
$ws.Range("H6").Value = "Symphony"
$ws.Range("J6").Value = 240
$ws.Range("D7").Value = 44377
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 20000
$ws.Range("L7").Value = 21000
$ws.Range("M7").Value = 20333
$ws.Range("P7").Value = 508
$ws.Range("D8").Value = 44377
$ws.Range("J8").Value = 60
$ws.Range("N8").Value = "$/caja 40 unidades"
$ws.Range("P8").Value = 538
$ws.Range("Q8").Value = 40
$ws.Range("D9").Value = 44435
$ws.Range("H9").Value = "Madrigal"
$ws.Range("J9").Value = 160
$ws.Range("K9").Value = 19000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 19500
$ws.Range("N9").Value = "$/caja 40 unidades"
$ws.Range("P9").Value = 488
$ws.Range("Q9").Value = 40
$ws.Range("D10").Value = 44391
$ws.Range("J10").Value = 140
$ws.Range("K10").Value = 21000
$ws.Range("L10").Value = 22000
$ws.Range("M10").Value = 21500
$ws.Range("P10").Value = 538
$ws.Range("D11").Value = 44419
$ws.Range("H11").Value = "Symphony"
$ws.Range("J11").Value = 150
$ws.Range("N11").Value = "$/caja 50 unidades"
$ws.Range("P11").Value = 430
$ws.Range("Q11").Value = 50
$ws.Range("D12").Value = 44405
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 21000
$ws.Range("L12").Value = 22000
$ws.Range("M12").Value = 21500
$ws.Range("N12").Value = "$/caja 40 unidades"
$ws.Range("P12").Value = 538
$ws.Range("Q12").Value = 40
$ws.Range("D13").Value = 44426
$ws.Range("H13").Value = "Madrigal"
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 20000
$ws.Range("M13").Value = 19500
$ws.Range("P13").Value = 488
$ws.Range("D14").Value = 44370
$ws.Range("H14").Value = "Argentina(o)"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 21000
$ws.Range("M14").Value = 20429
$ws.Range("N14").Value = "$/caja 50 unidades"
$ws.Range("P14").Value = 409
$ws.Range("Q14").Value = 50
$ws.Range("D15").Value = 44370
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 22000
$ws.Range("L15").Value = 23000
$ws.Range("M15").Value = 22500
$ws.Range("P15").Value = 562
$ws.Range("D16").Value = 44167
$ws.Range("H16").Value = "Española"
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 13000
$ws.Range("L16").Value = 14000
$ws.Range("M16").Value = 13500
$ws.Range("N16").Value = "$/caja 30 unidades"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 450
$ws.Range("Q16").Value = 30
$ws.Range("D17").Value = 44433
$ws.Range("J17").Value = 160
$ws.Range("K17").Value = 19000
$ws.Range("L17").Value = 20000
$ws.Range("M17").Value = 19500
$ws.Range("P17").Value = 488
$ws.Range("D18").Value = 44363
$ws.Range("D19").Value = 44356
$ws.Range("H19").Value = "Argentina(o)"
$ws.Range("J19").Value = 120
$ws.Range("K19").Value = 19000
$ws.Range("L19").Value = 20000
$ws.Range("M19").Value = 19500
$ws.Range("N19").Value = "$/caja 50 unidades"
$ws.Range("P19").Value = 390
$ws.Range("Q19").Value = 50
$ws.Range("D20").Value = 44160
$ws.Range("H20").Value = "Madrigal"
$ws.Range("K20").Value = 14000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 14500
$ws.Range("N20").Value = "$/caja 40 unidades"
$ws.Range("O20").Value = "Región de Coquimbo"
$ws.Range("P20").Value = 362
$ws.Range("Q20").Value = 40
